$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "61.269.94"
$ws.Range("E2").Value = "  -3.16%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.957.32"
$ws.Range("E3").Value = "  -2.21%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.11%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.90"
$ws.Range("E5").Value = "  -4.05%  "

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.61"
$ws.Range("E6").Value = "  -5.76%  "

# Row 7: USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.00%  "

# Row 8: XRP
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.558"
$ws.Range("E8").Value = "  -0.06%  "

# Row 9: LidoStakedEther
$ws.Range("D9").Value = "2.964.95"
$ws.Range("E9").Value = "  -2.16%  "

# Row 10: Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.111"
$ws.Range("E10").Value = "  -1.75%  "

# Row 11: Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.05"
$ws.Range("E11").Value = "  -5.60%  "

# Row 12: Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.363"
$ws.Range("E12").Value = "  -1.14%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.479.50"
$ws.Range("E13").Value = "  -1.97%  "

# Row 14: TRON
$ws.Range("E14").Value = "  +1.58%  "

# Row 15: WrappedBTC
$ws.Range("D15").Value = "61.350.91"
$ws.Range("E15").Value = "  -2.95%  "

# Row 16: Avalanche
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.51"
$ws.Range("E16").Value = "  -2.51%  "

# Row 17: WrappedEther
$ws.Range("D17").Value = "2.969.22"
$ws.Range("E17").Value = "  -1.68%  "

# Row 18: ShibaInu
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000144"
$ws.Range("E18").Value = "  -4.46%  "

# Row 19: Polkadot
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.08"
$ws.Range("E19").Value = "  -0.60%  "

# Row 20: Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.86"
$ws.Range("E20").Value = "  -1.76%  "

# Row 21: BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "371.05"
$ws.Range("E21").Value = "  -6.55%  "

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.63"
$ws.Range("E22").Value = "  -1.29%  "

# Row 23: Dai
$ws.Range("E23").Value = "  -0.09%  "

# Row 24: Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.25"
$ws.Range("E24").Value = "  +0.13%  "

# Row 25: WrappedeETH
$ws.Range("D25").Value = "3.099.43"
$ws.Range("E25").Value = "  -1.58%  "

# Row 26: Polygon
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.464"
$ws.Range("E26").Value = "  -0.60%  "

# Row 27: Kaspa
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.186"
$ws.Range("E27").Value = "  -1.78%  "

# Row 28: Binance-PegBSC-USD
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  +0.36%  "

# Row 29: PEPE
$ws.Range("D29").Value = "0.0₃0900"
$ws.Range("E29").Value = "  -7.66%  "

# Row 30: InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.05"
$ws.Range("E30").Value = "  -8.11%  "

# Row 31: USDe
$ws.Range("E31").Value = "  +0.02%  "

# Row 32: PancakeSwap
$ws.Range("E32").Value = "  -3.99%  "

# Row 33: EthereumClassic
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.14"
$ws.Range("E33").Value = "  -1.76%  "

# Row 34: Monero
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "159.57"
$ws.Range("E34").Value = "  -1.81%  "

# Row 35: NEARProtocol
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.49"
$ws.Range("E35").Value = "  -5.56%  "

# Row 36: Aptos
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.79"
$ws.Range("E36").Value = "  -4.13%  "

# Row 37: Fetch.AI
$ws.Range("E37").Value = "  -7.79%  "

# Row 38: ImmutableX
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.24"
$ws.Range("E38").Value = "  -5.67%  "

# Row 39: Stacks
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.51"
$ws.Range("E39").Value = "  -5.94%  "

# Row 40: Maker -> OKB
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.11"
$ws.Range("E40").Value = "  -1.54%  "

# Row 41: OKB -> Maker
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.391.87"
$ws.Range("E41").Value = "  -6.06%  "

# Row 42: Filecoin
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.81"
$ws.Range("E42").Value = "  -3.38%  "

# Row 43: Mantle
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.663"
$ws.Range("E43").Value = "  -0.48%  "

# Row 44: EnergySwap
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.58"
$ws.Range("E44").Value = "  -5.43%  "

# Row 45: Hedera
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0582"
$ws.Range("E45").Value = "  -3.57%  "

# Row 46: RenderToken
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.11"
$ws.Range("E46").Value = "  +0.32%  "

# Row 47: FirstDigitalUSD
$ws.Range("E47").Value = "  +0.15%  "

# Row 48: VeChain
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0240"
$ws.Range("E48").Value = "  -4.32%  "

# Row 49: Bittensor -> Stellar
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0939"
$ws.Range("E49").Value = "  -1.12%  "

# Row 50: Stellar -> WhiteBITCoin
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.40"

# Row 51: WhiteBITCoin -> Bittensor
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "262.37"
$ws.Range("E51").Value = "  -3.13%  "

